$wb = $excel.ActiveWorkbook

# --- Add the new "Employees" sheet as the first sheet ---
$wsEmployees = $wb.Worksheets.Add()
$wsEmployees.Name = "Employees"
$wsEmployees.Move($wb.Worksheets.Item(1))

# Header row
$wsEmployees.Range("A1").Value = "ID"
$wsEmployees.Range("B1").Value = "Name"
$wsEmployees.Range("C1").Value = "Department"

# Data rows
$wsEmployees.Range("A2").Value = 1
$wsEmployees.Range("B2").Value = "Alice"
$wsEmployees.Range("C2").Value = "HR"

$wsEmployees.Range("A3").Value = 2
$wsEmployees.Range("B3").Value = "Bob"
$wsEmployees.Range("C3").Value = "Legal"

$wsEmployees.Range("A4").Value = 3
$wsEmployees.Range("B4").Value = "Charlie"
$wsEmployees.Range("C4").Value = "IT"

$wsEmployees.Range("A5").Value = 4
$wsEmployees.Range("B5").Value = "Diana"
$wsEmployees.Range("C5").Value = "Marketing"

$wsEmployees.Range("A6").Value = 4
$wsEmployees.Range("B6").Value = "Dion"
$wsEmployees.Range("C6").Value = "Marketing"

# --- Update the existing "Login" sheet ---
$wsLogin = $wb.Worksheets.Item("Login")

# Remove the trailing empty rows (4:6)
$wsLogin.Rows.Item(4).Resize(3).Delete()

# Update credentials shown in the example data
$wsLogin.Range("B2").Value = "password1"
$wsLogin.Range("A3").Value = "sbuda@gmail.com"

# --- Make Login the active sheet/tab with the same selection as the source file ---
$wsLogin.Activate() | Out-Null
$wsLogin.Range("D6").Select() | Out-Null
